$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate Polish header labels to English (commit: "Polish names -> English names")
$ws.Range("B1").Value = "Publication DOI"
$ws.Range("C1").Value = "Number of compound in publication"
$ws.Range("E1:E121").Value = "Activity [nM]"

# Reflect the edited range as the active selection, like the author left it
$ws.Range("E2:E121").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
